$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 98; $r++) {
    $src = $ws.Range("AJ$r")
    $dst = $ws.Range("AK$r")
    $src.Cut($dst)
}
